# "add master data and python utility"
# The language master-data sheet is trimmed down to only the "eng" and
# "hin" rows (French/Arabic/Kannada/Tamil rows are removed), and the
# remaining data rows lose their inherited header-style formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the fra/ara/kan/tam rows (rows 4 through 7), leaving header +
# eng + (what becomes) hin.
$ws.Rows("4:7").Delete()

# The surviving data rows (2:3) no longer carry the bordered/centered
# header style - reset them back to the default "Normal" style.
$ws.Range("A2:D3").Style = "Normal"

# Row 3 now holds the Hindi entry instead of French.
$ws.Range("A3").Value = "hin"
$ws.Range("B3").Value = "हिन्दी"
$ws.Range("C3").Value = "भारोपीय"
$ws.Range("D3").Value = "Hindi"

# Match the author's last on-screen selection.
$ws.Range("F2").Select()
